$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values for the new rows (17-25) ---
$ws.Range("A17").Value = "SCRIPT/P02P01A/um2103.ssb"
$ws.Range("A18").Value = "SCRIPT/T01P01A/um2202.ssb"
$ws.Range("B18").Value = 69
$ws.Range("C18").Value = " We came running when we heard\nthe world\'s in danger!"
$ws.Range("D18").Value = " Как только мы узнали, что мир\nв опасности, мы сразу же прибежали сюда!"
$ws.Range("E18").Value = " Ëàë óïìûëï íú ôèîàìé, œóï íéñ\nâ ïðàòîïòóé, íú òñàèô çå ðñéáåçàìé òýäà!"
$ws.Range("B19").Value = 72
$ws.Range("C19").Value = " I tell you, this is no time to be\nsoaking in the Hot Spring!"
$ws.Range("D19").Value = " Говорю вам, нет времени мокнуть\nв Горячих Источниках!"
$ws.Range("E19").Value = " Ãïâïñý âàí, îåó âñåíåîé íïëîôóû\nâ Ãïñÿœéö Éòóïœîéëàö!"
$ws.Range("B20").Value = 75
$ws.Range("C20").Value = " So...[K]we\'ve come but…"
$ws.Range("D20").Value = " Так...[K] Мы здесь, но..."
$ws.Range("E20").Value = " Óàë...[K] Íú èäåòû, îï..."
$ws.Range("B21").Value = 78
$ws.Range("C21").Value = " What are we supposed to\ndo now...?"
$ws.Range("D21").Value = " Что нам теперь делать?.."
$ws.Range("E21").Value = " Œóï îàí óåðåñý äåìàóý?.."
$ws.Range("A22").Value = "SCRIPT/P02P01A/um2403.ssb"
$ws.Range("B22").Value = 41
$ws.Range("C22").Value = " It turns out we weren\'t being\nhelpful in any way…"
$ws.Range("D22").Value = " Оказалось, что мы не можем\nничем помочь..."
$ws.Range("E22").Value = " Ïëàèàìïòý, œóï íú îå íïçåí\nîéœåí ðïíïœý..."
$ws.Range("A23").Value = "SCRIPT/P02P01A/um2503.ssb"
$ws.Range("B23").Value = 44
$ws.Range("C23").Value = " So we came back to the\nHot Spring."
$ws.Range("D23").Value = " Поэтому мы вернулись на Горячие\nИсточники."
$ws.Range("E23").Value = " Ðïüóïíô íú âåñîôìéòý îà Ãïñÿœéå\nÉòóïœîéëé."
$ws.Range("B24").Value = 47
$ws.Range("C24").Value = " I feel guilty for being relaxed\nas I say this, but...[K]we\'re counting on you!"
$ws.Range("D24").Value = " Я немного стыжусь того, что тут\nотдыхаю, но...[K] Мы надеемся на вас!"
$ws.Range("E24").Value = " Ÿ îåíîïãï òóúçôòý óïãï, œóï óôó\nïóäúöàý, îï...[K] Íú îàäååíòÿ îà âàò!"
$ws.Range("B25").Value = 50
$ws.Range("C25").Value = " Stop the planet\'s paralysis or\nsomething! You can do it!"
$ws.Range("D25").Value = " Остановите планетарный паралич\nили что там ещё! Вы справитесь!"
$ws.Range("E25").Value = " Ïòóàîïâéóå ðìàîåóàñîúê ðàñàìéœ\néìé œóï óàí åþæ! Âú òðñàâéóåòý!"

# --- Borders for the separator rows (17 and 21): thin bottom border ---
$ws.Range("A17:B17").Borders.Item(9).LineStyle = 1
$ws.Range("C17:E17").Borders.Item(9).LineStyle = 1
$ws.Range("A21:B21").Borders.Item(9).LineStyle = 1
$ws.Range("C21:E21").Borders.Item(9).LineStyle = 1

# --- Row heights to match the authored layout ---
$ws.Rows.Item(17).RowHeight = 43.2
$ws.Rows.Item(18).RowHeight = 43.2
$ws.Rows.Item(19).RowHeight = 21.6
$ws.Rows.Item(20).RowHeight = 14.4
$ws.Rows.Item(21).RowHeight = 25.8
$ws.Rows.Item(22).RowHeight = 43.2
$ws.Rows.Item(23).RowHeight = 43.2
$ws.Rows.Item(24).RowHeight = 31.8
$ws.Rows.Item(25).RowHeight = 21.6

# --- View state: scrolled down with D23 selected ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("D23").Select()
